# Apply weekly update: insert two new price records (rows 1131-1132) for
# "Terminal Hortofrutícola Agro Chillán - Cebolla" and shift the existing
# rows down, keeping the trailing rows (now 1190-1191) as duplicates of
# the previous last two rows (as captured by the OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 1131; this pushes
# the existing data (previously in rows 1131-1189) down to rows 1133-1191,
# which already reproduces the target end-of-sheet rows 1190-1191.
$ws.Rows("1131:1132").Insert()

# --- New row 1131 ---
$ws.Range("A1131").Value = 7
$ws.Range("B1131").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C1131").Value = "Ñuble"
$ws.Range("D1131").Value = 45147
$ws.Range("E1131").Value = 16
$ws.Range("F1131").Value = 100112004
$ws.Range("G1131").Value = "Cebolla"
$ws.Range("H1131").Value = "Sin especificar"
$ws.Range("I1131").Value = "1a (guarda)"
$ws.Range("J1131").Value = 180
$ws.Range("K1131").Value = 15000
$ws.Range("L1131").Value = 15000
$ws.Range("M1131").Value = 15000
$ws.Range("N1131").Value = "`$/malla 25 kilos"
$ws.Range("O1131").Value = "Región del Maule"
$ws.Range("P1131").Value = 600
$ws.Range("Q1131").Value = 25
$ws.Range("R1131").Value = "Hortaliza"

# --- New row 1132 ---
$ws.Range("A1132").Value = 7
$ws.Range("B1132").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C1132").Value = "Ñuble"
$ws.Range("D1132").Value = 45147
$ws.Range("E1132").Value = 16
$ws.Range("F1132").Value = 100112004
$ws.Range("G1132").Value = "Cebolla"
$ws.Range("H1132").Value = "Sin especificar"
$ws.Range("I1132").Value = "2a (guarda)"
$ws.Range("J1132").Value = 150
$ws.Range("K1132").Value = 13000
$ws.Range("L1132").Value = 13000
$ws.Range("M1132").Value = 13000
$ws.Range("N1132").Value = "`$/malla 25 kilos"
$ws.Range("O1132").Value = "Región del Maule"
$ws.Range("P1132").Value = 520
$ws.Range("Q1132").Value = 25
$ws.Range("R1132").Value = "Hortaliza"
